# Update the "base de datos" (worker mora data) on sheet Hoja1.
# The two existing detail rows (16 and 17) had their "Periodo Mora"
# (column E, shared string "1810"/"1811") and "Valor Mora" (column F)
# values swapped between them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Periodo Mora -> "1811", Valor Mora -> 40000
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 40000

# Row 17: Periodo Mora -> "1810", Valor Mora -> 20000
$ws.Range("E17").Value = "1810"
$ws.Range("F17").Value = 20000
